$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Helper: find the 1-based index of the paragraph whose text contains
# the given substring.
# ------------------------------------------------------------------
function Get-ParagraphIndex($needle) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

# ------------------------------------------------------------------
# 1) "Лабораторная работа № 2" -> "Лабораторная работа № 4"
# ------------------------------------------------------------------
$labIdx = Get-ParagraphIndex("Лабораторная работа")
$labPara = $d.Paragraphs.Item($labIdx)
$labPara.Range.Find.Execute("2", $false, $false, $false, $false, $false, $true, 1, $false, "4", 2)

# ------------------------------------------------------------------
# 2) Remove the existing "_GoBack" bookmark (it will be re-inserted
#    later, right before the closing guillemet of the lab title).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 3) Drop the whole paragraph "«Численные методы одномерной
#    минимизации" - its content is being replaced by new title text
#    that now lives in the following paragraph.
# ------------------------------------------------------------------
$oldTitleIdx = Get-ParagraphIndex("Численные методы одномерной минимизации")
$d.Paragraphs.Item($oldTitleIdx).Range.Delete()

# ------------------------------------------------------------------
# 4) Turn the remaining "с использованием производной»" paragraph
#    into "«Численные методы многомерной оптимизации с
#    использованием производных первого порядка»", keeping the
#    trailing guillemet run intact and re-inserting the "_GoBack"
#    bookmark right before it.
# ------------------------------------------------------------------
$tailIdx = Get-ParagraphIndex("с использованием производной»")
$tailPara = $d.Paragraphs.Item($tailIdx)
$tailStart = $tailPara.Range.Start

$oldRunLen = "с использованием производной".Length
$oldRun = $d.Range($tailStart, $tailStart + $oldRunLen)
$newTitle = "«Численные методы многомерной оптимизации с использованием производных первого порядка"
$oldRun.Text = $newTitle

$bmPos = $tailStart + $newTitle.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
